$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): B1 and D1 get the "text" number format (49 = "@")
# applied first so the new header-style xf (font2/fill2/border2 + numFmt49)
# is allocated before the data-row style, matching the target style index order.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "@"

# --- Data rows 2-3: convert the national-id / mobile numbers stored as
# numbers into left-zero-padded text values, with the "text" number format.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"

$ws.Range("B2").Value = "0020123448"
$ws.Range("B3").Value = "0020123447"
$ws.Range("D2").Value = "09123456789"
$ws.Range("D3").Value = "09123456789"

# --- Remaining (empty) data rows 4-20 in columns B and D also switch to the
# text format so newly typed IDs keep their leading zeros.
$ws.Range("B4:B20").NumberFormat = "@"
$ws.Range("D4:D20").NumberFormat = "@"

# Rows 2 and 3 shrink back to the default row height once the taller typed
# content is gone.
$ws.Rows("2:3").AutoFit()

# --- Column layout: columns B and D were previously grouped together with
# column C under one width entry; give them their own column entries (width
# stays visually the same).
$ws.Columns("B").ColumnWidth = 45.5
$ws.Columns("D").ColumnWidth = 45.5

# --- Selection / view: the sheet now opens scrolled to the default position
# with D2 selected (instead of being scrolled to D1 with E4 selected).
$ws.Range("D2").Select()
